$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay as plain text (matching source data format)
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.737.25'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.466.41'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '578.58'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.15'
$ws.Range('E6').Value = '  +3.76%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +1.61%  '
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.401'
$ws.Range('E11').Value = '  +4.14%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '4.060.21'
$ws.Range('E12').Value = '  +2.27%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '29.83'
$ws.Range('E13').Value = '  +5.47%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.129'
$ws.Range('E14').Value = '  +2.40%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.457.61'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '62.781.12'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('E18').Value = '  +3.25%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '14.35'
$ws.Range('E19').Value = '  +5.41%  '
$ws.Range('E20').Value = '  +2.60%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '388.13'
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('E22').Value = '  +2.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '74.51'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.610.81'
$ws.Range('E25').Value = '  +2.43%  '
$ws.Range('E26').Value = '  +1.53%  '
$ws.Range('E27').Value = '  -9.98%  '
$ws.Range('E28').Value = '  +2.90%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.17'
$ws.Range('E30').Value = '  +2.22%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.14'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  -1.77%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '23.71'
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.29'
$ws.Range('E35').Value = '  +5.24%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.09'
$ws.Range('E36').Value = '  +2.69%  '
$ws.Range('E37').Value = '  +7.70%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '31.44'
$ws.Range('E38').Value = '  +22.34%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '169.44'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.503.46'
$ws.Range('E40').Value = '  +2.32%  '
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.798'
$ws.Range('E42').Value = '  +2.49%  '
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '42.36'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.48'
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.72'
$ws.Range('E45').Value = '  +4.43%  '
$ws.Range('E46').Value = '  +3.70%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.609.67'
$ws.Range('E47').Value = '  +5.94%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '23.12'
$ws.Range('E48').Value = '  +1.64%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.22'
$ws.Range('E49').Value = '  +10.70%  '
$ws.Range('E50').Value = '  +1.28%  '
$ws.Range('E51').Value = '  +0.07%  '
